$d = $word.ActiveDocument

$d.Content.Find.Execute("{date}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2019-04-25", 2)

$d.Content.Find.Execute("{expediteur}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "رئيس مصلحة كتابة الضبط بالمحكمة الادارية بأكادير", 2)

$d.Content.Find.Execute("{destinataire}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "مديرية الموارد البشرية", 2)

$d.Content.Find.Execute("{num_order}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2019/00022", 2)

$d.Content.Find.Execute("{order}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "01", 2)

$d.Content.Find.Execute("{text}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "شهادة التسليم تتعلق  بلوائح الخبرات الغير منجزة برسم سنتي 2017 و 2018", 2)

$d.Content.Find.Execute("{nb_copy}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "02", 2)

$d.Content.Find.Execute("{remarque}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "  نرجعها  لكم بعد القيام بالمطلوب، تبعا لإرسالكم عدد 424/2019 بتاريخ 01/04/2019، ", 2)
